# repull data, push all data, mean calculation
# Update the "dSF" (column F) values for rows 2-5 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -5
$ws.Range("F3").Value = -2
$ws.Range("F4").Value = -1
$ws.Range("F5").Value = 6
